$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "71.140.86"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +3.16%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.608.24"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +2.52%  "
$ws.Range("E4").Value = "  -0.01%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "603.62"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.78%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "174.82"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.81%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "3.604.11"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +2.60%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.618"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.05%  "
$ws.Range("E9").Value = "  -0.03%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.203"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +7.54%  "
$ws.Range("E11").Value = "  +8.29%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.594"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +2.23%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "47.35"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.45%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.0000281"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +2.25%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "4.180.99"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +2.53%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "8.50"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.17%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "620.98"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.83%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "3.604.83"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +2.46%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "71.151.54"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +3.12%  "
$ws.Range("E20").Value = "  -1.75%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "17.58"
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.894"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.73%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "9.33"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -16.39%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "16.19"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.71%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "98.04"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.04%  "
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("E27").Value = "  +1.97%  "
$ws.Range("E28").Value = "  +0.01%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "34.33"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +4.95%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "9.37"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +1.20%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "8.57"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.28%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.12"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.79%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "7.32"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +5.67%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.32"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.89%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "631.24"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -1.36%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "3.78"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +8.73%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("E38").Value = "  +1.75%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.0487"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +7.29%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "57.72"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.88%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.10%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.144"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +6.90%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "3.409.70"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.62%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.329"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("D45").Value = "0.0₃0724"
$ws.Range("E45").Value = "  +3.78%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.73"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +7.20%  "
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "3.01"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +9.70%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "33.20"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("E49").Value = "  +0.97%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "132.88"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.18%  "
